$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.ClearFormats()
}

# Row 37/38: NEARProtocol and Aptos swap positions (with updated price/volume)
Set-TextCell 'B37' 'Aptos'
Set-TextCell 'C37' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D37' '7.92'
Set-TextCell 'E37' '  +0.82%  '
Set-TextCell 'B38' 'NEARProtocol'
Set-TextCell 'C38' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D38' '5.71'
Set-TextCell 'E38' '  -6.08%  '

# Price / Volume(1h) updates
Set-TextCell 'D2' '66.874.19'
Set-TextCell 'E2' '  +0.19%  '
Set-TextCell 'D3' '3.485.96'
Set-TextCell 'E3' '  +0.75%  '
Set-TextCell 'D4' '1.00'
Set-TextCell 'E4' '  -0.13%  '
Set-TextCell 'D5' '605.99'
Set-TextCell 'E5' '  +1.22%  '
Set-TextCell 'D6' '145.14'
Set-TextCell 'E6' '  -1.35%  '
Set-TextCell 'D7' '3.484.66'
Set-TextCell 'E7' '  +0.73%  '
Set-TextCell 'D8' '1.00'
Set-TextCell 'E8' '  -0.02%  '
Set-TextCell 'D9' '0.476'
Set-TextCell 'E9' '  -1.29%  '
Set-TextCell 'D10' '0.140'
Set-TextCell 'E10' '  -0.49%  '
Set-TextCell 'E11' '  +7.35%  '
Set-TextCell 'D12' '0.417'
Set-TextCell 'E12' '  -1.50%  '
Set-TextCell 'D13' '0.0000212'
Set-TextCell 'E13' '  +0.26%  '
Set-TextCell 'D14' '4.084.75'
Set-TextCell 'E14' '  +0.57%  '
Set-TextCell 'D15' '31.06'
Set-TextCell 'E15' '  -1.45%  '
Set-TextCell 'D16' '3.493.85'
Set-TextCell 'E16' '  +0.21%  '
Set-TextCell 'D17' '66.538.61'
Set-TextCell 'E17' '  -0.70%  '
Set-TextCell 'E18' '  +0.12%  '
Set-TextCell 'D19' '10.79'
Set-TextCell 'E19' '  +9.19%  '
Set-TextCell 'D20' '6.27'
Set-TextCell 'E20' '  -2.09%  '
Set-TextCell 'D21' '15.31'
Set-TextCell 'E21' '  +0.77%  '
Set-TextCell 'D22' '427.32'
Set-TextCell 'E22' '  -2.38%  '
Set-TextCell 'D23' '0.601'
Set-TextCell 'E23' '  -2.57%  '
Set-TextCell 'D24' '78.97'
Set-TextCell 'E24' '  +0.41%  '
Set-TextCell 'E25' '  +0.09%  '
Set-TextCell 'D26' '3.626.61'
Set-TextCell 'E26' '  +0.45%  '
Set-TextCell 'E27' '  -1.52%  '
Set-TextCell 'D28' '9.69'
Set-TextCell 'E28' '  -1.05%  '
Set-TextCell 'D29' '8.13'
Set-TextCell 'E29' '  -2.13%  '
Set-TextCell 'D30' '2.50'
Set-TextCell 'E30' '  +1.22%  '
Set-TextCell 'D31' '1.54'
Set-TextCell 'E31' '  -3.33%  '
Set-TextCell 'D32' '0.999'
Set-TextCell 'E32' '  -0.34%  '
Set-TextCell 'D33' '0.165'
Set-TextCell 'E33' '  +1.32%  '
Set-TextCell 'D34' '25.28'
Set-TextCell 'E34' '  -0.07%  '
Set-TextCell 'D35' '1.77'
Set-TextCell 'E35' '  -1.69%  '
Set-TextCell 'E36' '  +0.01%  '
Set-TextCell 'D39' '1.00'
Set-TextCell 'E39' '  -0.09%  '
Set-TextCell 'D40' '174.91'
Set-TextCell 'E40' '  +0.49%  '
Set-TextCell 'D41' '0.0890'
Set-TextCell 'E41' '  +0.75%  '
Set-TextCell 'D42' '5.31'
Set-TextCell 'E42' '  -0.93%  '
Set-TextCell 'D43' '0.891'
Set-TextCell 'E43' '  +0.63%  '
Set-TextCell 'D44' '1.96'
Set-TextCell 'E44' '  -11.21%  '
Set-TextCell 'D45' '46.18'
Set-TextCell 'E45' '  -0.26%  '
Set-TextCell 'D46' '27.93'
Set-TextCell 'E46' '  -6.19%  '
Set-TextCell 'D47' '1.21'
Set-TextCell 'E47' '  -2.88%  '
Set-TextCell 'D48' '7.32'
Set-TextCell 'E48' '  -2.04%  '
Set-TextCell 'D49' '2.38'
Set-TextCell 'E49' '  -1.12%  '
Set-TextCell 'D50' '0.971'
Set-TextCell 'E50' '  -1.11%  '
Set-TextCell 'D51' '0.244'
Set-TextCell 'E51' '  -0.30%  '
